$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Sample table: style swap + header-row shading/text-color + bold labels
# ---------------------------------------------------------------------------
$tbl = $d.Tables(1)
$tbl.Style = "TableGrid"

# Header row (row 1): blue fill + white bold text for each of the 3 cells.
for ($c = 1; $c -le 3; $c++) {
    $cell = $tbl.Cell(1, $c)

    $cell.Shading.Texture = 0
    $cell.Shading.ForegroundPatternColor = -16777216
    $cell.Shading.BackgroundPatternColor = 13998939

    $full = $cell.Range
    $txt = $d.Range($full.Start, $full.End - 1)
    $txt.Font.Color = 16777215
}

# Data rows: bold the first-column label ("Text editing" / "Tables").
for ($r = 2; $r -le 3; $r++) {
    $cell = $tbl.Cell($r, 1)
    $full = $cell.Range
    $txt = $d.Range($full.Start, $full.End - 1)
    $txt.Bold = 1
}

# ---------------------------------------------------------------------------
# 2. Page header: centered, gray "DOCX JS Editor"
# ---------------------------------------------------------------------------
$sec = $d.Sections(1)

$hdr = $sec.Headers(1)
$hdr.Range.InsertAfter("DOCX JS Editor")
$hdr.Range.Paragraphs(1).Style = "Header"
$hdr.Range.ParagraphFormat.Alignment = 1
$hdrTxt = $hdr.Range.Duplicate
$hdrTxt.SetRange($hdr.Range.Start, $hdr.Range.Start + 14)
$hdrTxt.Font.Color = 8421504
$hdrTxt.Font.Size = 9
$hdrTxt.Font.SizeBi = 9

# ---------------------------------------------------------------------------
# 3. Page footer: centered, gray "Page 1"
# ---------------------------------------------------------------------------
$ftr = $sec.Footers(1)
$ftr.Range.InsertAfter("Page 1")
$ftr.Range.Paragraphs(1).Style = "Footer"
$ftr.Range.ParagraphFormat.Alignment = 1
$ftrTxt = $ftr.Range.Duplicate
$ftrTxt.SetRange($ftr.Range.Start, $ftr.Range.Start + 6)
$ftrTxt.Font.Color = 8421504
$ftrTxt.Font.Size = 9
$ftrTxt.Font.SizeBi = 9
